# Generate Report for Handback
#
# Updates the handback-status workbook with a fresh handoff/handback
# report for the "1c3734cb-541c-4ddd-ad8f-20d7fda2ecab.md" source file:
#   - Overview sheet: bump "Latest HO Xliff Generate Date"
#   - zh-cn sheet: bump "Correspond Handoff Datetime" / "Correspond Handback DateTime"
#   - de-de sheet: bump "Correspond Handoff Datetime" / "Correspond Handback DateTime"

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-09-03 10:52:57"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-09-03 10:52:53"
$zhcn.Range("K2").Value = "2016-09-03 10:53:15"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-09-03 10:52:57"
$dede.Range("K2").Value = "2016-09-03 10:53:21"
